$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 453.6
$ws.Range("I19").Value = 430.30768
$ws.Range("J19").Value = 467.36365
$ws.Range("K19").Value = 430.30768
$ws.Range("L19").Value = 467.36365
$ws.Range("M19").Value = -255.30768
$ws.Range("N19").Value = -817.36365
$ws.Range("H28").Value = 530223.5
$ws.Range("I28").Value = 855416.4
$ws.Range("J28").Value = 1785
$ws.Range("K28").Value = 855416.4
$ws.Range("L28").Value = 1785
$ws.Range("M28").Value = -854931.4
$ws.Range("N28").Value = -2755
$ws.Range("H107").Value = 654069.9
$ws.Range("I107").Value = 694886.4
$ws.Range("K107").Value = 694886.4
$ws.Range("M107").Value = -692966.4
$ws.Range("H132").Value = 13294.927
$ws.Range("I132").Value = 14422.676
$ws.Range("J132").Value = 2863.25
$ws.Range("K132").Value = 43268.028
$ws.Range("L132").Value = 8589.75
$ws.Range("M132").Value = -40738.028
$ws.Range("N132").Value = -13649.75
$ws.Range("H138").Value = 6947091.5
$ws.Range("I138").Value = 1955561.9
$ws.Range("J138").Value = 10206866
$ws.Range("K138").Value = 5866685.699999999
$ws.Range("L138").Value = 30620598
$ws.Range("M138").Value = -5861545.699999999
$ws.Range("N138").Value = -30630878
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 60679.117
$ws.Range("I2").Value = 68569.664
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 68569.664
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -68456.664
$ws.Range("N2").Value = -1726
$ws.Range("H32").Value = 18955.574
$ws.Range("I32").Value = 4541.6855
$ws.Range("K32").Value = 4541.6855
$ws.Range("M32").Value = -4254.6855
$ws.Range("H43").Value = 10377
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10377
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10377
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -11003
$ws.Range("H45").Value = 959
$ws.Range("I45").Value = 958.36365
$ws.Range("J45").Value = 961.3333
$ws.Range("K45").Value = 958.36365
$ws.Range("L45").Value = 961.3333
$ws.Range("M45").Value = -581.36365
$ws.Range("N45").Value = -1715.3333
$ws.Range("H63").Value = 9978.75
$ws.Range("I63").Value = 11351.25
$ws.Range("J63").Value = 9292.5
$ws.Range("K63").Value = 11351.25
$ws.Range("L63").Value = 9292.5
$ws.Range("M63").Value = -10665.25
$ws.Range("N63").Value = -10664.5
$ws.Range("H66").Value = 9978.75
$ws.Range("I66").Value = 11351.25
$ws.Range("J66").Value = 9292.5
$ws.Range("K66").Value = 56756.25
$ws.Range("L66").Value = 46462.5
$ws.Range("M66").Value = -53324.25
$ws.Range("N66").Value = -53326.5
$ws.Range("H116").Value = 60679.117
$ws.Range("I116").Value = 68569.664
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 68569.664
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = -66275.664
$ws.Range("N116").Value = -6088
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 60679.117
$ws.Range("I3").Value = 68569.664
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 68569.664
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -68455.664
$ws.Range("N3").Value = -1728
$ws.Range("H20").Value = 1231.9642
$ws.Range("I20").Value = 955.1579
$ws.Range("J20").Value = 1816.3334
$ws.Range("K20").Value = 955.1579
$ws.Range("L20").Value = 1816.3334
$ws.Range("M20").Value = -708.1579
$ws.Range("N20").Value = -2310.3334
$ws.Range("H134").Value = 2846.1277
$ws.Range("I134").Value = 2078.1765
$ws.Range("J134").Value = 4854.615
$ws.Range("K134").Value = 6234.529500000001
$ws.Range("L134").Value = 14563.845
$ws.Range("M134").Value = -3699.529500000001
$ws.Range("N134").Value = -19633.845
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 968.25
$ws.Range("I105").Value = 1017.1429
$ws.Range("J105").Value = 899.8
$ws.Range("K105").Value = 1017.1429
$ws.Range("L105").Value = 899.8
$ws.Range("M105").Value = 729.8570999999999
$ws.Range("N105").Value = -4393.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 14465.777
$ws.Range("H71").Value = 14465.777
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 5000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -5576
$ws.Range("H80").Value = 2540.9092
$ws.Range("I80").Value = 2615.3845
$ws.Range("J80").Value = 2433.3333
$ws.Range("K80").Value = 2615.3845
$ws.Range("L80").Value = 2433.3333
$ws.Range("M80").Value = -1617.3845
$ws.Range("N80").Value = -4429.3333
$ws.Range("H81").Value = 5000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 5000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -6996
$ws.Range("H83").Value = 2540.9092
$ws.Range("I83").Value = 2615.3845
$ws.Range("J83").Value = 2433.3333
$ws.Range("K83").Value = 13076.9225
$ws.Range("L83").Value = 12166.6665
$ws.Range("M83").Value = -8084.922500000001
$ws.Range("N83").Value = -22150.6665
$ws.Range("H84").Value = 5000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 15000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -24984
$ws.Range("H113").Value = 1674.75
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 1666.3334
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 1666.3334
$ws.Range("M113").Value = 470
$ws.Range("N113").Value = -6006.3334
$ws.Range("H125").Value = 35900
$ws.Range("J125").Value = 35900
$ws.Range("L125").Value = 35900
$ws.Range("N125").Value = -40820
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3305.7778
$ws.Range("I7").Value = 2168
$ws.Range("J7").Value = 3533.3333
$ws.Range("K7").Value = 2168
$ws.Range("L7").Value = 3533.3333
$ws.Range("M7").Value = -2056
$ws.Range("N7").Value = -3757.3333
$ws.Range("H126").Value = 3305.7778
$ws.Range("I126").Value = 2168
$ws.Range("J126").Value = 3533.3333
$ws.Range("K126").Value = 6504
$ws.Range("L126").Value = 10599.9999
$ws.Range("M126").Value = -4034
$ws.Range("N126").Value = -15539.9999
$ws.Range("H132").Value = 3176.0732
$ws.Range("I132").Value = 2148.3103
$ws.Range("J132").Value = 5659.8335
$ws.Range("K132").Value = 6444.9309
$ws.Range("L132").Value = 16979.5005
$ws.Range("M132").Value = -3914.9309
$ws.Range("N132").Value = -22039.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1370.9445
$ws.Range("I122").Value = 1147
$ws.Range("J122").Value = 1879.909
$ws.Range("K122").Value = 3441
$ws.Range("L122").Value = 5639.727000000001
$ws.Range("M122").Value = -991
$ws.Range("N122").Value = -10539.727
$ws.Range("H136").Value = 3068.4167
$ws.Range("I136").Value = 1390.5555
$ws.Range("J136").Value = 5225.6665
$ws.Range("K136").Value = 4171.666499999999
$ws.Range("L136").Value = 15676.9995
$ws.Range("M136").Value = -1621.666499999999
$ws.Range("N136").Value = -20776.9995
